# lab03/Libro1.xlsx - "Se modifico el inventario del lab03"
#
# The PC DELL ALL IN ONE CORE I9 units now include 32GB of RAM, so the
# description column needs updating. All three inventory rows (l001-l003)
# shared the same description, so update all three the same way - this
# keeps them sharing a single entry in the workbook's shared-string table,
# exactly like the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDescription = "PC DELL ALL IN ONE CORE I9 mem 32gb"

$ws.Range("C2").Value = $newDescription
$ws.Range("C3").Value = $newDescription
$ws.Range("C4").Value = $newDescription

# The longer text no longer fits the old column width - resize column C
# (descripcion) to fit the new content.
$ws.Columns.Item(3).AutoFit()

# Leave the selection where the author left it before saving.
$ws.Range("C6").Select()
